$d = $word.ActiveDocument
$d.Content.Find.Execute("2024-02-28 Wednesday", $false, $false, $false, $false, $false, $true, 1, $false, "2024-02-29 Thursday", 2) | Out-Null
$d.Content.Find.Execute("926÷6=154, 2", $false, $false, $false, $false, $false, $true, 1, $false, "142÷9=15, 7", 2) | Out-Null
$d.Content.Find.Execute("277÷8=34, 5", $false, $false, $false, $false, $false, $true, 1, $false, "398÷7=56, 6", 2) | Out-Null
$d.Content.Find.Execute("899÷2=449, 1", $false, $false, $false, $false, $false, $true, 1, $false, "141÷5=28, 1", 2) | Out-Null
$d.Content.Find.Execute("549÷5=109, 4", $false, $false, $false, $false, $false, $true, 1, $false, "804÷8=100, 4", 2) | Out-Null
$d.Content.Find.Execute("597÷2=298, 1", $false, $false, $false, $false, $false, $true, 1, $false, "884÷2=442, 0", 2) | Out-Null
$d.Content.Find.Execute("657÷2=328, 1", $false, $false, $false, $false, $false, $true, 1, $false, "587÷6=97, 5", 2) | Out-Null
$d.Content.Find.Execute("990÷9=110, 0", $false, $false, $false, $false, $false, $true, 1, $false, "182÷3=60, 2", 2) | Out-Null
$d.Content.Find.Execute("579÷2=289, 1", $false, $false, $false, $false, $false, $true, 1, $false, "413÷9=45, 8", 2) | Out-Null
$d.Content.Find.Execute("640÷6=106, 4", $false, $false, $false, $false, $false, $true, 1, $false, "921÷4=230, 1", 2) | Out-Null
$d.Content.Find.Execute("668÷9=74, 2", $false, $false, $false, $false, $false, $true, 1, $false, "226÷5=45, 1", 2) | Out-Null
$d.Content.Find.Execute("594÷7=84, 6", $false, $false, $false, $false, $false, $true, 1, $false, "318÷5=63, 3", 2) | Out-Null
$d.Content.Find.Execute("702÷8=87, 6", $false, $false, $false, $false, $false, $true, 1, $false, "949÷8=118, 5", 2) | Out-Null
$d.Content.Find.Execute("856÷8=107, 0", $false, $false, $false, $false, $false, $true, 1, $false, "914÷8=114, 2", 2) | Out-Null
$d.Content.Find.Execute("655÷3=218, 1", $false, $false, $false, $false, $false, $true, 1, $false, "275÷3=91, 2", 2) | Out-Null
$d.Content.Find.Execute("750÷9=83, 3", $false, $false, $false, $false, $false, $true, 1, $false, "221÷3=73, 2", 2) | Out-Null
$d.Content.Find.Execute("390÷6=65, 0", $false, $false, $false, $false, $false, $true, 1, $false, "895÷8=111, 7", 2) | Out-Null
$d.Content.Find.Execute("799÷7=114, 1", $false, $false, $false, $false, $false, $true, 1, $false, "304÷7=43, 3", 2) | Out-Null
$d.Content.Find.Execute("708÷3=236, 0", $false, $false, $false, $false, $false, $true, 1, $false, "704÷8=88, 0", 2) | Out-Null
$d.Content.Find.Execute("569÷9=63, 2", $false, $false, $false, $false, $false, $true, 1, $false, "444÷7=63, 3", 2) | Out-Null
$d.Content.Find.Execute("221÷9=24, 5", $false, $false, $false, $false, $false, $true, 1, $false, "129÷6=21, 3", 2) | Out-Null
$d.Content.Find.Execute("330÷2=165, 0", $false, $false, $false, $false, $false, $true, 1, $false, "529÷5=105, 4", 2) | Out-Null
$d.Content.Find.Execute("110÷7=15, 5", $false, $false, $false, $false, $false, $true, 1, $false, "753÷9=83, 6", 2) | Out-Null
$d.Content.Find.Execute("747÷2=373, 1", $false, $false, $false, $false, $false, $true, 1, $false, "802÷2=401, 0", 2) | Out-Null
$d.Content.Find.Execute("318÷7=45, 3", $false, $false, $false, $false, $false, $true, 1, $false, "601÷6=100, 1", 2) | Out-Null
$d.Content.Find.Execute("624÷5=124, 4", $false, $false, $false, $false, $false, $true, 1, $false, "533÷4=133, 1", 2) | Out-Null
